# Avance de la secuencia de las piezas en la estacion
#
# Rewrites three existing log-entry paragraphs so that the English/technical
# terms inside them (Listener, thread, Runnable, Estacion, JPanel, Gui) are
# wrapped with <w:proofErr w:type="spellStart"/>...<w:proofErr w:type="spellEnd"/>
# markers (splitting the paragraph into several runs), fixes "esta" -> "está"
# in the Listener paragraph, and appends two brand-new log paragraphs at the
# end of the document.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Paragraph: "Agregada clase Listener, implementa thread ..." ----------
$fragListener = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t xml:space="preserve">Agregada clase </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Listener</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, implementa </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>thread</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> y es para escuchar cuando la cortadora </w:t></w:r>' +
    '<w:r><w:t>está</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> ocupada</w:t></w:r>' +
    '</w:p>'

# --- Paragraph: "Agregada la implementación de Runnable a la clase Estacion"
$fragRunnable = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t xml:space="preserve">Agregada la implementación de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Runnable</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> a la clase </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Estacion</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

# --- Paragraph: "Agregado un JPanel Estaciones panel ..." -----------------
$fragJPanel = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t xml:space="preserve">Agregado un </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>JPanel</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Estaciones panel para realizar pruebas con una </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Gui</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> respecto a las estaciones, de momento solo muestra una representación de la estación cortadora y un botón que la activa. Mientras la estación esta libre se muestra en verde y mientras esté ocupada en rojo</w:t></w:r>' +
    '</w:p>'

# --- Two brand-new paragraphs appended at the end of the body -------------
$fragNew = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t>Agregada la implementación de recibir piezas a la estación</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t>La estación avanza la etapa de las piezas cuando las procesa</w:t></w:r>' +
    '</w:p>'

# Locate each paragraph by its distinctive original text (more robust than a
# hard-coded paragraph index) and replace its content in place.
$lastParagraph = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -like "Agregada clase Listener*") {
        $p.Range.InsertXML($fragListener)
    }
    elseif ($t -like "Agregada la implementación de Runnable*") {
        $p.Range.InsertXML($fragRunnable)
    }
    elseif ($t -like "Agregado un JPanel*") {
        $p.Range.InsertXML($fragJPanel)
    }

    $lastParagraph = $p
}

# Append the two new paragraphs right after the very last paragraph in the
# document body (the "... estación dobladora" one), before the sectPr.
$endRange = $d.Range($lastParagraph.Range.End, $lastParagraph.Range.End)
$endRange.InsertXML($fragNew)
